# Add vendor rows to the "articels" sheet: rows 10-13 in column A get the
# text "/addven" (replacing what used to be a long run of empty placeholder
# rows down to row 99).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articels")

for ($r = 10; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = "/addven"
}
